$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'compression knee'
$ws.Cells.Item(2, 1).Value = 'hockey pads men'
$ws.Cells.Item(3, 1).Value = 'black basketball'
$ws.Cells.Item(4, 1).Value = 'leg pads baseball'
$ws.Cells.Item(5, 1).Value = 'running tights men'
$ws.Cells.Item(6, 1).Value = 'snowboarding pads'
$ws.Cells.Item(7, 1).Value = 'compression leggings youth boys'
$ws.Cells.Item(8, 1).Value = 'knee pads wrestling men'
$ws.Cells.Item(9, 1).Value = 'workout pants for girls'
$ws.Cells.Item(10, 1).Value = 'basketball team gear'
$ws.Cells.Item(11, 1).Value = 'basketball apparel for girls'
$ws.Cells.Item(12, 1).Value = 'boys basketball compression leggings'
$ws.Cells.Item(13, 1).Value = 'compressions pants'
$ws.Cells.Item(14, 1).Value = 'bjj tights men'
$ws.Cells.Item(15, 1).Value = 'compression basketball pants'
$ws.Cells.Item(16, 1).Value = 'mens baseball pants'
$ws.Cells.Item(17, 1).Value = 'best compression leggings'
$ws.Cells.Item(18, 1).Value = 'softball pants black'
$ws.Cells.Item(19, 1).Value = 'sliding pants youth'
$ws.Cells.Item(20, 1).Value = 'basketballs pants'
$ws.Cells.Item(21, 1).Value = 'sport knee pads'
$ws.Cells.Item(22, 1).Value = 'men s knee pads'
$ws.Cells.Item(23, 1).Value = 'volleyball knee pads black'
$ws.Cells.Item(24, 1).Value = 'pack of leggings'
$ws.Cells.Item(25, 1).Value = 'cycling knee protector'
$ws.Cells.Item(26, 1).Value = 'knee pads sleeve basketball'
$ws.Cells.Item(27, 1).Value = 'compression tights recovery'
$ws.Cells.Item(28, 1).Value = 'volleyball pads'
$ws.Cells.Item(29, 1).Value = 'padded tights for basketball'
$ws.Cells.Item(30, 1).Value = 'small knee pads'
$ws.Cells.Item(31, 1).Value = 'soccer apparel youth'
$ws.Cells.Item(32, 1).Value = 'cheap knee pads'
$ws.Cells.Item(33, 1).Value = 'baseball gear boys'
$ws.Cells.Item(34, 1).Value = 'knee pads six six one'
$ws.Cells.Item(35, 1).Value = 'compression pad'
$ws.Cells.Item(36, 1).Value = 'free volleyball'
$ws.Cells.Item(37, 1).Value = 'gym pad'
$ws.Cells.Item(38, 1).Value = 'compression pants for boys basketball'
$ws.Cells.Item(39, 1).Value = 'girl knee pads'
$ws.Cells.Item(40, 1).Value = 'boys leggings sports youth'
$ws.Cells.Item(41, 1).Value = 'youth football pants with pads small'
$ws.Cells.Item(42, 1).Value = 'mens wrestling knee pads'
$ws.Cells.Item(43, 1).Value = 'sport compression pants boys'
$ws.Cells.Item(44, 1).Value = 'softball pants for girls'
$ws.Cells.Item(45, 1).Value = 'compression leggings'
$ws.Cells.Item(46, 1).Value = 'mens lacrosse pads'
$ws.Cells.Item(47, 1).Value = 'softball pants youth girls black'
$ws.Cells.Item(48, 1).Value = 'sport pants men'
$ws.Cells.Item(49, 1).Value = 'thigh compression tights'
$ws.Cells.Item(50, 1).Value = 'youth knee compression'
$ws.Cells.Item(51, 1).Value = 'basketball sleeve knee pads'
$ws.Cells.Item(52, 1).Value = 'leggings boys basketball'
$ws.Cells.Item(53, 1).Value = 'basketball gear'
$ws.Cells.Item(54, 1).Value = 'lacrosse youth pads'
$ws.Cells.Item(55, 1).Value = 'running tights for men'
$ws.Cells.Item(56, 1).Value = 'compressions leggings for men'
$ws.Cells.Item(57, 1).Value = 'knee gel pads'
$ws.Cells.Item(58, 1).Value = 'mens big and tall pants'
$ws.Cells.Item(59, 1).Value = 'basketball protection'
$ws.Cells.Item(60, 1).Value = 'knee braces for men xxl'
$ws.Cells.Item(61, 1).Value = 'leggings with mesh girls'
$ws.Cells.Item(62, 1).Value = 'little black pants guaranteed to fit'
$ws.Cells.Item(63, 1).Value = '3/4 pants men'
$ws.Cells.Item(64, 1).Value = 'mens capri tights'
$ws.Cells.Item(65, 1).Value = 'boys basketball'
$ws.Cells.Item(66, 1).Value = 'girls volleyball spandex'
$ws.Cells.Item(67, 1).Value = 'girls softball pants'
$ws.Cells.Item(68, 1).Value = 'lcl knee support'
$ws.Cells.Item(69, 1).Value = 'arthritis equipment'
$ws.Cells.Item(70, 1).Value = 'volleyball knee pads'
$ws.Cells.Item(71, 1).Value = 'running compression men'
$ws.Cells.Item(72, 1).Value = 'baseball pants youth xl'
$ws.Cells.Item(73, 1).Value = 'youth xs knee pads basketball'
$ws.Cells.Item(74, 1).Value = 'under pants for basketball'
$ws.Cells.Item(75, 1).Value = 'guy legging'
$ws.Cells.Item(76, 1).Value = 'under armor warm pants'
$ws.Cells.Item(77, 1).Value = 'spandex leggings men thermal'
$ws.Cells.Item(78, 1).Value = 'mems thermal leggings'
$ws.Cells.Item(79, 1).Value = 'men compression 3 4 pants'
$ws.Cells.Item(80, 1).Value = 'men compression leggings 3 4'
$ws.Cells.Item(81, 1).Value = 'men compression pants adidas'
$ws.Cells.Item(82, 1).Value = 'men compression tights 3 4'
$ws.Cells.Item(83, 1).Value = 'mens compression 3 4 tights'
$ws.Cells.Item(84, 1).Value = 'mens compression 3 4 leggings'
$ws.Cells.Item(85, 1).Value = 'eastbay baseball pants'
$ws.Cells.Item(86, 1).Value = 'eastbay leggings'
$ws.Cells.Item(87, 1).Value = 'eastbay tights'
$ws.Cells.Item(88, 1).Value = 'mcdavid knee pads youth'
$ws.Cells.Item(89, 1).Value = 'basketball leggings nike'
$ws.Cells.Item(90, 1).Value = 'nba basketball pants'
$ws.Cells.Item(91, 1).Value = 'basketball 3 4 leggings'
$ws.Cells.Item(92, 1).Value = 'elbow knee wrist pads for youth'
$ws.Cells.Item(93, 1).Value = 'telsa thermals men'
$ws.Cells.Item(94, 1).Value = 'wintergear compression leggings men'
$ws.Cells.Item(95, 1).Value = 'dry skin tights'
$ws.Cells.Item(96, 1).Value = 'mens workoit tights'
$ws.Cells.Item(97, 1).Value = 'track leggings'
$ws.Cells.Item(98, 1).Value = 'track tights for men'
$ws.Cells.Item(99, 1).Value = 'tesla wintergear compression shirt'
$ws.Cells.Item(100, 1).Value = 'basketball apparel youth'
